$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# "Nombre d'heures par personne" report - fill in the hours worked by
# Mathieu GAILLARD (row 8) across the seven project phases (cols E:K).
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 5
$ws.Range("G8").Value = 18
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 5
$ws.Range("J8").Value = 6
$ws.Range("K8").Value = 1

# The per-person (L) and per-phase (row 13) totals are driven by existing
# SUM formulas already in the sheet, so they recalculate automatically.

# Reflect where the user ended up looking afterwards: scrolled slightly
# right and focused on K9.
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("K9").Select()
